$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

# Duplicate the formatting/content of the existing last row (date style on
# col A, shared-text cells on col G/H which both carry the same "2"/"KK.MI"
# values needed for the new rows too) onto the two new rows, then overwrite
# the numeric values that actually differ.
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy($ws.Range("A" + $newRow1 + ":H" + $newRow1))
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy($ws.Range("A" + $newRow2 + ":H" + $newRow2))

# Row 100
$ws.Cells.Item($newRow1, 1).Value = 45457.2916666667
$ws.Cells.Item($newRow1, 2).Value = 0
$ws.Cells.Item($newRow1, 3).Value = 2
$ws.Cells.Item($newRow1, 4).Value = 2
$ws.Cells.Item($newRow1, 5).Value = 2
$ws.Cells.Item($newRow1, 6).Value = 2
# Columns G ("2") and H ("KK.MI") already match after the row copy above.

# Row 101
$ws.Cells.Item($newRow2, 1).Value = 45460.3122453704
$ws.Cells.Item($newRow2, 2).Value = 5400
$ws.Cells.Item($newRow2, 3).Value = 2
$ws.Cells.Item($newRow2, 4).Value = 2
$ws.Cells.Item($newRow2, 5).Value = 2
$ws.Cells.Item($newRow2, 6).Value = 2
# Columns G ("2") and H ("KK.MI") already match after the row copy above.
